# VejlederMatrix.xlsx - rebuild the "Resultat" occupancy matrix.
# The old sheet only tracked one pairwise collision per row (column C, with a
# distinct "Optaget (X og Y)" label per combination). The new layout spreads
# the occupied markers across many columns (one "Optaget" marker per actual
# collision, all sharing the same generic "Optaget" text) and widens the used
# range out to T9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultat")

# Keep columns A/B as-is (names), but size A:C like the authored sheet.
$ws.Columns("A").ColumnWidth = 9.666666666666666
$ws.Columns("B").ColumnWidth = 9.666666666666666
$ws.Columns("C").ColumnWidth = 8.333333333333332

# Row 2: previously only C2 ("Optaget (MKR og DRK)") -> now H2, L2
$ws.Range("C2").ClearContents()
$ws.Range("H2").Value = "Optaget"
$ws.Range("L2").Value = "Optaget"

# Row 3: previously only C3 ("Optaget (DRK og LOD)") -> now C3,D3,F3,J3,O3
$ws.Range("C3").Value = "Optaget"
$ws.Range("D3").Value = "Optaget"
$ws.Range("F3").Value = "Optaget"
$ws.Range("J3").Value = "Optaget"
$ws.Range("O3").Value = "Optaget"

# Row 4: previously only C4 ("Optaget (UOP og MKR)") -> now E4,I4,M4
$ws.Range("C4").ClearContents()
$ws.Range("E4").Value = "Optaget"
$ws.Range("I4").Value = "Optaget"
$ws.Range("M4").Value = "Optaget"

# Row 5: previously only C5 ("Optaget (LOD og SOD)") -> now E5,G5,I5,K5,N5,Q5
$ws.Range("C5").ClearContents()
$ws.Range("E5").Value = "Optaget"
$ws.Range("G5").Value = "Optaget"
$ws.Range("I5").Value = "Optaget"
$ws.Range("K5").Value = "Optaget"
$ws.Range("N5").Value = "Optaget"
$ws.Range("Q5").Value = "Optaget"

# Row 6: previously only C6 ("Optaget (SOD og UOP)") -> now C6,D6,F6,H6,J6,L6,O6,R6
$ws.Range("C6").Value = "Optaget"
$ws.Range("D6").Value = "Optaget"
$ws.Range("F6").Value = "Optaget"
$ws.Range("H6").Value = "Optaget"
$ws.Range("J6").Value = "Optaget"
$ws.Range("L6").Value = "Optaget"
$ws.Range("O6").Value = "Optaget"
$ws.Range("R6").Value = "Optaget"

# Row 7: previously only C7 ("Optaget (SOD og MKR)") -> now S7
$ws.Range("C7").ClearContents()
$ws.Range("S7").Value = "Optaget"

# Row 8: previously only C8 ("Optaget (DRK og MKR)") -> now G8,K8,N8
$ws.Range("C8").ClearContents()
$ws.Range("G8").Value = "Optaget"
$ws.Range("K8").Value = "Optaget"
$ws.Range("N8").Value = "Optaget"

# Row 9: previously only C9 ("Optaget (LOD og SOD)") -> now M9,P9,T9
$ws.Range("C9").ClearContents()
$ws.Range("M9").Value = "Optaget"
$ws.Range("P9").Value = "Optaget"
$ws.Range("T9").Value = "Optaget"

# Match the authored selection (cursor left on U8 after filling the matrix).
$ws.Range("U8").Select()
